$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Construction Materials"
$ws.Range("B2").Value = "1"

# Row 3
$ws.Range("A3").Value = "Labor — Construction"
$ws.Range("B3").Value = "1"

# Row 4
$ws.Range("A4").Value = " Vehicles"
$ws.Range("B4").Value = "1"

# Row 5
$ws.Range("A5").Value = "Labor"
$ws.Range("B5").Value = "11"

# Row 6
$ws.Range("A6").Value = "Transformers"
$ws.Range("B6").Value = "13"

# Row 7
$ws.Range("A7").Value = "Electrical Components"
$ws.Range("B7").Value = "6"

# Row 8
$ws.Range("A8").Value = "Heating, Ventilation  Air Conditioning"
$ws.Range("B8").Value = "HVAC) Equipment (2"
